$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.01565089521678647
$ws.Cells.Item(2, 3).Value = 0.01484250298558193
$ws.Cells.Item(3, 2).Value = 0.001292249644803861
$ws.Cells.Item(3, 3).Value = 0.0009980755223185536
$ws.Cells.Item(4, 2).Value = 0.008562355981155206
$ws.Cells.Item(4, 3).Value = 0.008250474730512054
$ws.Cells.Item(5, 2).Value = 0.01173737137161623
$ws.Cells.Item(5, 3).Value = 0.01097223691847164
$ws.Cells.Item(6, 2).Value = 0.01225909768248579
$ws.Cells.Item(6, 3).Value = 0.01179365268582781
$ws.Cells.Item(7, 2).Value = 0.007676332646453125
$ws.Cells.Item(7, 3).Value = 0.007327655420201334
$ws.Cells.Item(8, 2).Value = 0.01027795651280915
$ws.Cells.Item(8, 3).Value = 0.009744268181040247
$ws.Cells.Item(9, 2).Value = 0.00577035119431221
$ws.Cells.Item(9, 3).Value = 0.005206448214256969
$ws.Cells.Item(10, 2).Value = 0.005953404137153302
$ws.Cells.Item(10, 3).Value = 0.00562710768736027
$ws.Cells.Item(11, 2).Value = 0.01490371428760107
$ws.Cells.Item(11, 3).Value = 0.01427715904756715
$ws.Cells.Item(12, 2).Value = 0.001804160764515152
$ws.Cells.Item(12, 3).Value = 0.001615578798598139
$ws.Cells.Item(13, 2).Value = 0.00887606859305887
$ws.Cells.Item(13, 3).Value = 0.008537345648789741
$ws.Cells.Item(14, 2).Value = 0.009997324532678572
$ws.Cells.Item(14, 3).Value = 0.00941981122403138
$ws.Cells.Item(15, 2).Value = 0.01249441337709281
$ws.Cells.Item(15, 3).Value = 0.01145037798942304
$ws.Cells.Item(16, 2).Value = 0.01075147302666425
$ws.Cells.Item(16, 3).Value = 0.01009987996601927
$ws.Cells.Item(17, 2).Value = 0.006116094199110228
$ws.Cells.Item(17, 3).Value = 0.005558892711175062
$ws.Cells.Item(18, 2).Value = 0.01031164216388764
$ws.Cells.Item(18, 3).Value = 0.009658649065965561
$ws.Cells.Item(19, 2).Value = 0.003831036454417572
$ws.Cells.Item(19, 3).Value = 0.003337130333700376
$ws.Cells.Item(20, 2).Value = 0.009936649004475872
$ws.Cells.Item(20, 3).Value = 0.009448718377917883
$ws.Cells.Item(21, 2).Value = 0.004479622553025265
$ws.Cells.Item(21, 3).Value = 0.004207748601397533
$ws.Cells.Item(22, 2).Value = 0.01209097048881631
$ws.Cells.Item(22, 3).Value = 0.01161376583036185
$ws.Cells.Item(23, 2).Value = 0.009017807743265677
$ws.Cells.Item(23, 3).Value = 0.008210484095044639
$ws.Cells.Item(24, 2).Value = 0.01010725658484811
$ws.Cells.Item(24, 3).Value = 0.00942132768668579
$ws.Cells.Item(25, 2).Value = 0.01222294875677461
$ws.Cells.Item(25, 3).Value = 0.01134305818350381
$ws.Cells.Item(26, 2).Value = 0.008452894945820428
$ws.Cells.Item(26, 3).Value = 0.007860420642530291
$ws.Cells.Item(27, 2).Value = 0.002585792365291421
$ws.Cells.Item(27, 3).Value = 0.002224053545289363
$ws.Cells.Item(28, 2).Value = 0.007608985076355341
$ws.Cells.Item(28, 3).Value = 0.007148314766172914
$ws.Cells.Item(29, 2).Value = 0.002431686678464668
$ws.Cells.Item(29, 3).Value = 0.002111994023015221
$ws.Cells.Item(30, 2).Value = 0.01655684758039638
$ws.Cells.Item(30, 3).Value = 0.01602284539581034
$ws.Cells.Item(31, 2).Value = 0.008819521322976821
$ws.Cells.Item(31, 3).Value = 0.008325469885120826
$ws.Cells.Item(32, 2).Value = 0.005507916550793147
$ws.Cells.Item(32, 3).Value = 0.005155766096489033
$ws.Cells.Item(33, 2).Value = 0.00918637286084714
$ws.Cells.Item(33, 3).Value = 0.008720767823559057
$ws.Cells.Item(34, 2).Value = 0.02218293317661213
$ws.Cells.Item(34, 3).Value = 0.02126390096473138
$ws.Cells.Item(35, 2).Value = 0.003891417411931921
$ws.Cells.Item(35, 3).Value = 0.003568812625069351
$ws.Cells.Item(36, 2).Value = 0.01061061810269276
$ws.Cells.Item(36, 3).Value = 0.009774345747532776
$ws.Cells.Item(37, 2).Value = 0.008104646224635385
$ws.Cells.Item(37, 3).Value = 0.007424518947062246
$ws.Cells.Item(38, 2).Value = 0.01069549817705062
$ws.Cells.Item(38, 3).Value = 0.01028141932814774
$ws.Cells.Item(39, 2).Value = 0.01127317596540012
$ws.Cells.Item(39, 3).Value = 0.01056947778967031
$ws.Cells.Item(40, 2).Value = 0.01355781914934823
$ws.Cells.Item(40, 3).Value = 0.01277478825732086
$ws.Cells.Item(41, 2).Value = 0.02046381824609615
$ws.Cells.Item(41, 3).Value = 0.01961376769917476
$ws.Cells.Item(42, 2).Value = 0.007222150649643902
$ws.Cells.Item(42, 3).Value = 0.006845687248900739
$ws.Cells.Item(43, 2).Value = 0.01374596070291992
$ws.Cells.Item(43, 3).Value = 0.01290581148818935
$ws.Cells.Item(44, 2).Value = 0.008924804062935089
$ws.Cells.Item(44, 3).Value = 0.008434898876428395
$ws.Cells.Item(45, 2).Value = 0.007747021260089229
$ws.Cells.Item(45, 3).Value = 0.007112784288590504
$ws.Cells.Item(46, 2).Value = 0.01149269492089652
$ws.Cells.Item(46, 3).Value = 0.01094730445112662
$ws.Cells.Item(47, 2).Value = 0.008683449328431361
$ws.Cells.Item(47, 3).Value = 0.008215098048639376
$ws.Cells.Item(48, 2).Value = 0.005761256039537013
$ws.Cells.Item(48, 3).Value = 0.005373034089407951
$ws.Cells.Item(49, 2).Value = 0.005720603205221458
$ws.Cells.Item(49, 3).Value = 0.004999894448816802
$ws.Cells.Item(50, 2).Value = 0.01284252534826277
$ws.Cells.Item(50, 3).Value = 0.01199056054463199
$ws.Cells.Item(51, 2).Value = 0.001598413865435728
$ws.Cells.Item(51, 3).Value = 0.001048529232748658
